$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the PICTURE table row (row 27): the PICTURE_DATA/mediumblob/picture-data
# column is replaced by a PICTURE_PATH/varchar/picture-address column.
$ws.Range("B27").Value = "PICTURE_PATH"
$ws.Range("C27").Value = "varchar"
$ws.Range("D27").Value = "图片地址"

# Move the active selection to D51 (was C52).
$ws.Range("D51").Select()
